$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - metric_2D / PC
$ws.Range("D2").Value = 0.0693
$ws.Range("F2").Value = 0.1089

# Row 3 - metric_breast_cancer / PC
$ws.Range("C3").Value = 0.0808
$ws.Range("D3").Value = 0.0439

# Row 4 - metric_load_iris / PC
$ws.Range("C4").Value = 0.0533
$ws.Range("D4").Value = 0.0067
$ws.Range("F4").Value = 0.08

# Row 5 - metric_load_wine / PC
$ws.Range("C5").Value = 0.0506
$ws.Range("D5").Value = 0.0169
$ws.Range("F5").Value = 0.0337

# Row 6 - indices_PC_LabelCorrection_before_fix_OCPC
$ws.Range("G6").Value = 0.06845
$ws.Range("H6").Value = 0.0002817725
$ws.Range("I6").Value = 0.01678608054311667

# Row 8 - indices_PC_LabelCorrection_after_fix_OCPC
$ws.Range("J8").Value = 0.0228
$ws.Range("K8").Value = 0.00065686
$ws.Range("L8").Value = 0.02562928013035091

# Row 9 - indices_CL_after_fix_OCPC
$ws.Range("J9").Value = 0.04355000000000001
$ws.Range("K9").Value = 0.001585495833333333
$ws.Range("L9").Value = 0.03981828516314249

# Row 10 - metric_2D / LOF
$ws.Range("C10").Value = 0.0446
$ws.Range("D10").Value = 0.0396
$ws.Range("F10").Value = 0.1089

# Row 11 - metric_breast_cancer / LOF
$ws.Range("C11").Value = 0.0861
$ws.Range("D11").Value = 0.0281
$ws.Range("F11").Value = 0.0457

# Row 12 - metric_load_iris / LOF
$ws.Range("C12").Value = 0.1067
$ws.Range("D12").Value = 0.0267
$ws.Range("F12").Value = 0.06

# Row 13 - metric_load_wine / LOF
$ws.Range("C13").Value = 0.1292
$ws.Range("D13").Value = 0.0281
$ws.Range("F13").Value = 0.0506

# Row 14 - indices_PC_LabelCorrection_before_fix_LOF
$ws.Range("G14").Value = 0.09165000000000001
$ws.Range("H14").Value = 0.0009702525000000002
$ws.Range("I14").Value = 0.03114887638422934

# Row 16 - indices_PC_LabelCorrection_after_fix_LOF
$ws.Range("J16").Value = 0.02041666666666667
$ws.Range("K16").Value = 0.0002265380555555556
$ws.Range("L16").Value = 0.0150511812013395

# Row 17 - indices_CL_after_fix_LOF
$ws.Range("J17").Value = 0.0442
$ws.Range("K17").Value = 0.001397703333333333
$ws.Range("L17").Value = 0.03738587077136673
